$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Comments" -> "Comments/Issues"
$ws.Range("C3").Value = "Comments/Issues"

# OSD: "Remove silly option thing (sidemenu)" comment gains "Options button up in top bar."
$ws.Range("C9").Value = "Fixed losing focus by using a window property, like Aeon Nox 5 does. Options button up in top bar. Dismiss by a back button."

# OSD transparent background colour nicer - now marked Done with a comment
$ws.Range("B18").Value = "Done"
$ws.Range("C18").Value = "Make it blue like WMC."

# Use full guide ... - add a comment about the keymap script
$ws.Range("C23").Value = "Ctrl-g for guide with support of a script (goes with optional keymap)"

# Guide roll up/down with mouse wheel - split combined comment, keep first half
$ws.Range("C31").Value = "Done using a keymap. "

# Guide needs up/down/left/right buttons - extend comment
$ws.Range("C32").Value = "Autoscroll might also be OK (but not so aggressive as the My Addons list!). Problem is, no button seems to be able to control the EPGGrid."

# Guide rows need to be a little bigger - add new comment about parameterising EPG sizes
$ws.Range("C33").Value = "Parameterise the EPG row size and font size (maybe even put in settings)"

# Update the saved view/selection to match the new scroll position
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
